$d = $word.ActiveDocument

# Fixed duplicated values when computing institutional factors.
# Updated all graphs and tables (p-values table in Fig 5).

# institutional: 0.01 -> 0.06
$d.Content.Find.Execute("0.01", $true, $false, $false, $false, $false,
                         $true, 1, $false, "0.06", 2)

# Cod: 0.38 -> 0.52
$d.Content.Find.Execute("0.38", $true, $false, $false, $false, $false,
                         $true, 1, $false, "0.52", 2)

# Hake: 0.03 -> 0.06
$d.Content.Find.Execute("0.03", $true, $false, $false, $false, $false,
                         $true, 1, $false, "0.06", 2)
